# Fruta / hortaliza, semanal
# A new weekly price record is inserted as row 10 on the "Poroto verde" sheet
# (Feria Lagunitas de Puerto Montt). All existing records from row 10 down
# shift one row lower (row 10 -> 11, ..., row 63 -> 64); the dimension grows
# from A1:R63 to A1:R64.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at position 10; Excel pushes rows 10:63 down to 11:64
# and copies the row formatting (e.g. the date style on column D) along.
$ws.Rows("10:10").Insert()

# Populate the newly inserted row 10 with the new weekly record.
$ws.Range("A10").Value2 = 4
$ws.Range("B10").Value2 = "Feria Lagunitas de Puerto Montt"
$ws.Range("C10").Value2 = "Los Lagos"
$ws.Range("D10").Value2 = 44635
$ws.Range("E10").Value2 = 10
$ws.Range("F10").Value2 = 100112031
$ws.Range("G10").Value2 = "Poroto verde"
$ws.Range("H10").Value2 = "Magnum"
$ws.Range("I10").Value2 = "Primera"
$ws.Range("J10").Value2 = 40
$ws.Range("K10").Value2 = 29000
$ws.Range("L10").Value2 = 29000
$ws.Range("M10").Value2 = 29000
$ws.Range("N10").Value2 = "$/saco 25 kilos"
$ws.Range("O10").Value2 = "Región Metropolitana"
$ws.Range("P10").Value2 = 1160
$ws.Range("Q10").Value2 = 25
$ws.Range("R10").Value2 = "Hortaliza"
